$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cell H1 should use the same formatting as the other header
# cells (e.g. G1): bold, bordered, centered/top-aligned.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save column values for rows 2-12
$values = @(0, 0, 1, 0, 0, 0, 0, 1, 1, 1, 0)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
